$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 7; $r -le 24; $r++) {
    $ws.Cells.Item($r, 3).Clear()
    $ws.Cells.Item($r, 2).Value = 17.9
    $ws.Cells.Item($r, 3).Value = 260.8
}
